$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data (D, E) under the existing header row
$ws.Range("A2").Value = "D"
$ws.Range("B2").Value = "E"

# Match the header-row formatting (bold, centered, yellow fill) used for A1
$rng = $ws.Range("A2:B2")
$rng.Font.Bold = $true
$rng.Interior.Color = 65535
$rng.HorizontalAlignment = -4108

# Update the active selection as recorded in the sheet view
[void]$ws.Range("H14").Select()
